$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 26582.21
$ws.Range("I8").Value = 111272.555
$ws.Range("K8").Value = 333817.665
$ws.Range("M8").Value = -333678.665
# Row 17
$ws.Range("H17").Value = 1461.3469
$ws.Range("J17").Value = 1487.3829
$ws.Range("L17").Value = 4462.1487
$ws.Range("N17").Value = -4798.1487
# Row 51
$ws.Range("I51").Value = 9999
$ws.Range("J51").Value = 9666.333000000001
$ws.Range("K51").Value = 9999
$ws.Range("L51").Value = 9666.333000000001
$ws.Range("M51").Value = -9515
$ws.Range("N51").Value = -10634.333
# Row 76
$ws.Range("H76").Value = 91014240
$ws.Range("I76").Value = 178200.83
$ws.Range("J76").Value = 200017470
$ws.Range("K76").Value = 178200.83
$ws.Range("L76").Value = 200017470
$ws.Range("M76").Value = -177885.83
$ws.Range("N76").Value = -200018100
# Row 79
$ws.Range("H79").Value = 91014240
$ws.Range("I79").Value = 178200.83
$ws.Range("J79").Value = 200017470
$ws.Range("K79").Value = 178200.83
$ws.Range("L79").Value = 200017470
$ws.Range("M79").Value = -177108.83
$ws.Range("N79").Value = -200019654
# Row 94
$ws.Range("H94").Value = 1282.8889
$ws.Range("I94").Value = 818.25
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 818.25
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -367.25
$ws.Range("N94").Value = -5902
# Row 98
$ws.Range("H98").Value = 2028.8966
$ws.Range("I98").Value = 1036.4348
$ws.Range("K98").Value = 1036.4348
$ws.Range("M98").Value = 461.5652
# Row 100
$ws.Range("H100").Value = 8676.416999999999
$ws.Range("I100").Value = 3014.5
$ws.Range("K100").Value = 3014.5
$ws.Range("M100").Value = -2473.5
# Row 112
$ws.Range("H112").Value = 2576.4482
$ws.Range("J112").Value = 2620.6072
$ws.Range("L112").Value = 7861.821599999999
$ws.Range("N112").Value = -10077.8216
# Row 122
$ws.Range("H122").Value = 2028.8966
$ws.Range("I122").Value = 1036.4348
$ws.Range("K122").Value = 3109.3044
$ws.Range("M122").Value = -659.3044

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3529.8262
$ws.Range("I45").Value = 3288.389
$ws.Range("K45").Value = 3288.389
$ws.Range("M45").Value = -2911.389
# Row 97
$ws.Range("H97").Value = 1129.6875
$ws.Range("I97").Value = 1137.0714
$ws.Range("K97").Value = 1137.0714
$ws.Range("M97").Value = -641.0714
# Row 102
$ws.Range("H102").Value = 3171.5293
$ws.Range("I102").Value = 3194.8
$ws.Range("K102").Value = 3194.8
$ws.Range("M102").Value = -1572.8
# Row 122
$ws.Range("H122").Value = 3615.238
$ws.Range("I122").Value = 2276.5715
$ws.Range("K122").Value = 6829.7145
$ws.Range("M122").Value = -4379.7145

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3771.0212
$ws.Range("I134").Value = 2483.1516
$ws.Range("J134").Value = 6806.7144
$ws.Range("K134").Value = 7449.4548
$ws.Range("L134").Value = 20420.1432
$ws.Range("M134").Value = -4914.4548
$ws.Range("N134").Value = -25490.1432

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 3327.4285
$ws.Range("I16").Value = 3215.5
$ws.Range("J16").Value = 3999
$ws.Range("K16").Value = 3215.5
$ws.Range("L16").Value = 3999
$ws.Range("M16").Value = -2928.5
$ws.Range("N16").Value = -4573
# Row 36
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
# Row 86
$ws.Range("H86").Value = 7287.2666
$ws.Range("J86").Value = 7626
$ws.Range("L86").Value = 7626
$ws.Range("N86").Value = -9872
# Row 89
$ws.Range("H89").Value = 7287.2666
$ws.Range("J89").Value = 7626
$ws.Range("L89").Value = 38130
$ws.Range("N89").Value = -49362
# Row 99
$ws.Range("H99").Value = 4771
$ws.Range("I99").Value = 3636.5715
$ws.Range("K99").Value = 3636.5715
$ws.Range("M99").Value = -2138.5715
# Row 113
$ws.Range("H113").Value = 3327.4285
$ws.Range("I113").Value = 3215.5
$ws.Range("J113").Value = 3999
$ws.Range("K113").Value = 3215.5
$ws.Range("L113").Value = 3999
$ws.Range("M113").Value = -1045.5
$ws.Range("N113").Value = -8339
# Row 126
$ws.Range("H126").Value = 4771
$ws.Range("I126").Value = 3636.5715
$ws.Range("K126").Value = 10909.7145
$ws.Range("M126").Value = -8439.7145
# Row 131
$ws.Range("H131").Value = 59649.5
$ws.Range("I131").Value = 39999
$ws.Range("K131").Value = 39999
$ws.Range("M131").Value = -34959
# Row 134
$ws.Range("H134").Value = 3512.7708
$ws.Range("I134").Value = 2539.9333
$ws.Range("K134").Value = 7619.7999
$ws.Range("M134").Value = -5084.7999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 26
$ws.Range("H26").Value = 499.33334
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
# Row 36
$ws.Range("H36").Value = 60941.57
$ws.Range("I36").Value = 379.75
$ws.Range("J36").Value = 141690.67
$ws.Range("K36").Value = 1139.25
$ws.Range("L36").Value = 425072.01
$ws.Range("M36").Value = -970.25
$ws.Range("N36").Value = -425410.01
# Row 92
$ws.Range("H92").Value = 5001501
$ws.Range("I92").Value = 10000002
$ws.Range("K92").Value = 30000006
$ws.Range("M92").Value = -29998758
# Row 97
$ws.Range("H97").Value = 1085.2
$ws.Range("J97").Value = 1030.75
$ws.Range("L97").Value = 3092.25
$ws.Range("N97").Value = -4084.25
# Row 101
$ws.Range("H101").Value = 13211.6
$ws.Range("J101").Value = 13211.6
$ws.Range("L101").Value = 39634.8
$ws.Range("N101").Value = -44502.8
# Row 122
$ws.Range("H122").Value = 112210.78
$ws.Range("J122").Value = 112210.78
$ws.Range("L122").Value = 1009897.02
$ws.Range("N122").Value = -1014797.02
# Row 128
$ws.Range("H128").Value = 174996.5
$ws.Range("I128").Value = 174996.5
$ws.Range("K128").Value = 524989.5
$ws.Range("M128").Value = -520009.5
# Row 137
$ws.Range("H137").Value = 1893
$ws.Range("I137").Value = 1571.6666
$ws.Range("J137").Value = 2375
$ws.Range("K137").Value = 4714.9998
$ws.Range("L137").Value = 7125
$ws.Range("M137").Value = 385.0002000000004
$ws.Range("N137").Value = -17325

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 43483116
$ws.Range("I70").Value = 4141.643
$ws.Range("J70").Value = 111117070
$ws.Range("K70").Value = 4141.643
$ws.Range("L70").Value = 111117070
$ws.Range("M70").Value = -3871.643
$ws.Range("N70").Value = -111117610
# Row 73
$ws.Range("H73").Value = 43483116
$ws.Range("I73").Value = 4141.643
$ws.Range("J73").Value = 111117070
$ws.Range("K73").Value = 4141.643
$ws.Range("L73").Value = 111117070
$ws.Range("M73").Value = -3205.643
$ws.Range("N73").Value = -111118942
# Row 102
$ws.Range("H102").Value = 1359.1875
$ws.Range("I102").Value = 1375.5333
$ws.Range("K102").Value = 1375.5333
$ws.Range("M102").Value = 246.4666999999999
# Row 122
$ws.Range("H122").Value = 14747.25
$ws.Range("I122").Value = 12990
$ws.Range("K122").Value = 38970
$ws.Range("M122").Value = -36520
# Row 132
$ws.Range("H132").Value = 288599.8
$ws.Range("I132").Value = 436829.44
$ws.Range("J132").Value = 4493
$ws.Range("K132").Value = 1310488.32
$ws.Range("L132").Value = 13479
$ws.Range("M132").Value = -1307958.32
$ws.Range("N132").Value = -18539

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 12
$ws.Range("H12").Value = 650
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 650
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 650
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -990
# Row 55
$ws.Range("H55").Value = 950.5599999999999
$ws.Range("I55").Value = 398.11765
$ws.Range("J55").Value = 2124.5
$ws.Range("K55").Value = 398.11765
$ws.Range("L55").Value = 2124.5
$ws.Range("M55").Value = -225.11765
$ws.Range("N55").Value = -2470.5
# Row 93
$ws.Range("H93").Value = 1853.1538
$ws.Range("I93").Value = 1600
$ws.Range("K93").Value = 1600
$ws.Range("M93").Value = -352
# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
# Row 132
$ws.Range("H132").Value = 4408.467
$ws.Range("I132").Value = 3760.1667
$ws.Range("K132").Value = 11280.5001
$ws.Range("M132").Value = -8750.500100000001
# Row 136
$ws.Range("H136").Value = 2758.92
$ws.Range("I136").Value = 1831.9
$ws.Range("J136").Value = 4149.45
$ws.Range("K136").Value = 5495.700000000001
$ws.Range("L136").Value = 12448.35
$ws.Range("M136").Value = -2945.700000000001
$ws.Range("N136").Value = -17548.35

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 6621
$ws.Range("I62").Value = 3901
$ws.Range("J62").Value = 8434.333000000001
$ws.Range("K62").Value = 3901
$ws.Range("L62").Value = 8434.333000000001
$ws.Range("M62").Value = -3277
$ws.Range("N62").Value = -9682.333000000001
# Row 65
$ws.Range("H65").Value = 6621
$ws.Range("I65").Value = 3901
$ws.Range("J65").Value = 8434.333000000001
$ws.Range("K65").Value = 19505
$ws.Range("L65").Value = 42171.665
$ws.Range("M65").Value = -16385
$ws.Range("N65").Value = -48411.665
# Row 107
$ws.Range("H107").Value = 100708
$ws.Range("I107").Value = 143329
$ws.Range("J107").Value = 1259
$ws.Range("K107").Value = 429987
$ws.Range("L107").Value = 3777
$ws.Range("M107").Value = -428067
$ws.Range("N107").Value = -7617
# Row 132
$ws.Range("H132").Value = 2061.4614
$ws.Range("I132").Value = 901
$ws.Range("J132").Value = 5542.846
$ws.Range("K132").Value = 2703
$ws.Range("L132").Value = 16628.538
$ws.Range("M132").Value = -173
$ws.Range("N132").Value = -21688.538

Write-Host "Applied all Jenova_Profits market data updates"